# Workbook / sheet references
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Daily Amazon")

# --- Add the 3 new "Gotenks" sales rows to the "Daily Amazon" sheet ------
# Copy the date formatting (yyyy-mm-dd) from the last existing row down
# onto the 3 new rows so the new date cells match the rest of column A.
$ws2.Range("A67").Copy()
$ws2.Range("A68:A70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 68; $r -le 70; $r++) {
    $ws2.Cells.Item($r, 1).Value = 44927            # Date  (2023-01-01)
    $ws2.Cells.Item($r, 2).Value = "Gotenks"         # Item
    $ws2.Cells.Item($r, 3).Value = 34.99             # Sold Price
    $ws2.Cells.Item($r, 4).Value = 12.11             # Fee
    $ws2.Cells.Item($r, 5).Value = 2                 # Item Cost
    $ws2.Cells.Item($r, 6).Value = 20.88             # Profit
}

# --- Leave a single-cell selection on "Sheet1" (it is no longer active) -
[void]$ws1.Range("D23").Select()

# --- Make "Daily Amazon" the active sheet/tab and set its selection -----
[void]$ws2.Activate()
[void]$ws2.Range("E74").Select()
